{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"676\u00d79=6084\", \"498\u00d79=4482\"],\n  [\"423\u00d73=1269\", \"751\u00d74=3004\"],\n  [\"816\u00d73=2448\", \"209\u00d74=836\"],\n  [\"483\u00d77=3381\", \"495\u00d73=1485\"],\n  [\"133\u00d78=1064\", \"356\u00d76=2136\"],\n  [\"720\u00d72=1440\", \"380\u00d72=760\"],\n  [\"521\u00d76=3126\", \"865\u00d77=6055\"],\n  [\"400\u00d72=800\", \"736\u00d73=2208\"],\n  [\"115\u00d73=345\", \"719\u00d79=6471\"],\n  [\"571\u00d79=5139\", \"948\u00d74=3792\"],\n  [\"640\u00d74=2560\", \"647\u00d74=2588\"],\n  [\"812\u00d75=4060\", \"320\u00d73=960\"],\n  [\"221\u00d79=1989\", \"591\u00d78=4728\"],\n  [\"586\u00d76=3516\", \"863\u00d74=3452\"],\n  [\"797\u00d72=1594\", \"535\u00d78=4280\"],\n  [\"316\u00d73=948\", \"561\u00d77=3927\"],\n  [\"806\u00d74=3224\", \"489\u00d78=3912\"],\n  [\"857\u00d72=1714\", \"469\u00d72=938\"],\n  [\"598\u00d75=2990\", \"457\u00d74=1828\"],\n  [\"944\u00d72=1888\", \"670\u00d72=1340\"],\n  [\"116\u00d72=232\", \"135\u00d78=1080\"],\n  [\"218\u00d78=1744\", \"106\u00d76=636\"],\n  [\"752\u00d79=6768\", \"780\u00d72=1560\"],\n  [\"939\u00d77=6573\", \"116\u00d75=580\"],\n  [\"666\u00d73=1998\", \"863\u00d72=1726\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"676\u00d79=6084\"\n$find.Replacement.Text = \"498\u00d79=4482\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"423\u00d73=1269\"\n$find.Replacement.Text = \"751\u00d74=3004\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"816\u00d73=2448\"\n$find.Replacement.Text = \"209\u00d74=836\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"483\u00d77=3381\"\n$find.Replacement.Text = \"495\u00d73=1485\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"133\u00d78=1064\"\n$find.Replacement.Text = \"356\u00d76=2136\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"720\u00d72=1440\"\n$find.Replacement.Text = \"380\u00d72=760\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"521\u00d76=3126\"\n$find.Replacement.Text = \"865\u00d77=6055\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"400\u00d72=800\"\n$find.Replacement.Text = \"736\u00d73=2208\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"115\u00d73=345\"\n$find.Replacement.Text = \"719\u00d79=6471\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"571\u00d79=5139\"\n$find.Replacement.Text = \"948\u00d74=3792\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"640\u00d74=2560\"\n$find.Replacement.Text = \"647\u00d74=2588\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"812\u00d75=4060\"\n$find.Replacement.Text = \"320\u00d73=960\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"221\u00d79=1989\"\n$find.Replacement.Text = \"591\u00d78=4728\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"586\u00d76=3516\"\n$find.Replacement.Text = \"863\u00d74=3452\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"797\u00d72=1594\"\n$find.Replacement.Text = \"535\u00d78=4280\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"316\u00d73=948\"\n$find.Replacement.Text = \"561\u00d77=3927\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"806\u00d74=3224\"\n$find.Replacement.Text = \"489\u00d78=3912\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"857\u00d72=1714\"\n$find.Replacement.Text = \"469\u00d72=938\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"598\u00d75=2990\"\n$find.Replacement.Text = \"457\u00d74=1828\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"944\u00d72=1888\"\n$find.Replacement.Text = \"670\u00d72=1340\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"116\u00d72=232\"\n$find.Replacement.Text = \"135\u00d78=1080\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"218\u00d78=1744\"\n$find.Replacement.Text = \"106\u00d76=636\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"752\u00d79=6768\"\n$find.Replacement.Text = \"780\u00d72=1560\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"939\u00d77=6573\"\n$find.Replacement.Text = \"116\u00d75=580\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"666\u00d73=1998\"\n$find.Replacement.Text = \"863\u00d72=1726\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
